# Applies the "Added computation edges underlying graph" edit:
#  - Reorders the data rows on the 'PREFIX vs PTD' sheet (rows 3-7), swapping
#    the mathoverflow/slashdot_reply rows and moving wiki_talk_ca up / pushing
#    enron_email and askubuntu down one row, with refreshed metric values.
#  - Updates the three 'PREFIX vs PTD'-scoped defined names so their ranges
#    keep tracking the same (now relocated) rows.
#  - Leaves the final selection/active sheet on 'PREFIX vs PTD' (cell C6),
#    matching the workbook's new activeTab / tabSelected state.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("PROXIES")
$ws3 = $wb.Worksheets.Item("PREFIX vs PTD")

# --- Rewrite the data block (rows 3-7) on 'PREFIX vs PTD' ------------------
# New row order/content (network name, time_ptd, spearman, ktau, wktau):
$rows = @(
    @{ Row = 3; Network = "slashdot_reply"; C = 1100.8770999999999; D = 0.10100000000000001; E = 0.88;                 F = 0.84;                 G = 0.96 },
    @{ Row = 4; Network = "mathoverflow";   C = 3534.1477;          D = 0.16489999999999999; E = 0.81;                 F = 0.71;                 G = 0.93 },
    @{ Row = 5; Network = "wiki_talk_ca";   C = 4813.1759000000002; D = 0.28060000000000002; E = 0.56000000000000005;  F = 0.55000000000000004;  G = 0.92 },
    @{ Row = 6; Network = "enron_email";    C = 8922.3081000000002; D = 0.28260000000000002; E = 0.79;                 F = 0.77;                 G = 0.96 },
    @{ Row = 7; Network = "askubuntu";      C = 37782.423000000003; D = 0.49509999999999998; E = 0.74;                 F = 0.67;                 G = 0.93 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws3.Range("B$row").Value = $r.Network
    $ws3.Range("C$row").Value = $r.C
    $ws3.Range("D$row").Value = $r.D
    $ws3.Range("E$row").Value = $r.E
    $ws3.Range("F$row").Value = $r.F
    $ws3.Range("G$row").Value = $r.G
}

# --- Keep the defined names in sync with the relocated rows ----------------
$names = $wb.Names

function Set-NameRefersTo($qualifiedName, $newRefersTo) {
    for ($i = 1; $i -le $names.Count; $i++) {
        $n = $names.Item($i)
        if ($n.Name -eq $qualifiedName) {
            $n.RefersTo = $newRefersTo
            return
        }
    }
}

Set-NameRefersTo "PREFIX vs PTD!_12_mathoverflow_not_onbra_4" "='PREFIX vs PTD'!`$B`$4:`$G`$4"
Set-NameRefersTo "PREFIX vs PTD!_12_mathoverflow_not_onbra_4_3" "='PREFIX vs PTD'!`$B`$3:`$G`$3"
Set-NameRefersTo "PREFIX vs PTD!_12_mathoverflow_not_onbra_4_4" "='PREFIX vs PTD'!`$B`$7:`$G`$7"

# --- Update view/selection state --------------------------------------------
# PROXIES: selection moves from L9 to L6.
$ws1.Activate()
$ws1.Range("L6").Select()

# 'PREFIX vs PTD' becomes the active sheet/tab, with C6 selected.
$ws3.Activate()
$ws3.Range("C6").Select()
